# Fix sales data: update existing rows 2-20 (CarId, brandId, sales_date)
# and append new rows 21-31 with additional sales records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("D_SNEHAMAH001", "C_AJAYTHO001", "3399541C-60E4-4C7E-A106-1DEC04C42D30", "MAH104", "S_001", 45713.59510416666),
    @("D_POOJAMAR001", "C_ANITACOI001", "1435503D-2763-44EA-8EE9-C1C8EDC556BC", "MAR109", "S_002", 45663.74355324074),
    @("D_SENTHILTOY001", "C_TANVICHE001", "FD84C8B1-00C0-4005-B94E-CB123B4B26F5", "TOY103", "S_003", 45758.68761574074),
    @("D_SURESHBMW001", "C_AJAYTHO001", "75DE4879-925F-430B-B03C-B34C51A0A539", "BMW102", "S_004", 45734.46260416666),
    @("D_ARUNTAT002", "C_SNEHATHO001", "77B23DBF-D020-493A-AC55-8161AEA035D3", "TAT106", "S_005", 45683.71097222222),
    @("D_POOJAMAR001", "C_MEENAKAN001", "A1D4538F-C5B4-4F86-8C5F-9AC4B01C0D99", "MAR109", "S_006", 45779.428125),
    @("D_SNEHAMAH001", "C_SUNITASAL001", "1DB67510-AD35-4DDD-AA89-AC96EEC26DDC", "MAH104", "S_007", 45807.71866898148),
    @("D_REVATHIHON001", "C_MEENAKAN001", "55323C0C-C991-4E2B-B249-C5039352C717", "HON107", "S_008", 45799.72204861111),
    @("D_ARUNTAT001", "C_AJAYTHO001", "2E59616C-4807-4771-8B73-A4AC6F487F1E", "TAT106", "S_009", 45680.38193287037),
    @("D_GOPINATHMER001", "C_PREETIVEL001", "C183551C-DDE9-4758-B9B7-FDC17E964E12", "MER105", "S_010", 45672.67266203704),
    @("D_PRIYAHYU001", "C_SWETHACOI002", "5A3F306A-9043-405D-A509-0179CA88DA7E", "HYU101", "S_011", 45801.55597222222),
    @("D_SENTHILTOY001", "C_KIRANTIR001", "9635C27A-5E74-42AF-A15A-D1E7677D8F75", "TOY103", "S_012", 45759.42482638889),
    @("D_MEERAMER001", "C_MEENAKAN001", "F4D299B7-E898-4CC0-BB43-B0270D250DAF", "MER105", "S_013", 45784.61175925926),
    @("D_LAKSHMIFOR001", "C_NITINVEL001", "7DCDD228-40BE-4039-B9F2-EED2704497E1", "FOR108", "S_014", 45718.52344907408),
    @("D_SANJAYHYU001", "C_ANITACOI001", "637A54FB-6A83-4C76-8DF0-FAD88B0C7E63", "HYU101", "S_015", 45694.4),
    @("D_VIJAYHON001", "C_PRIYATIR001", "5ACFDEFC-070B-4BD0-97D8-8D6D1AA5E6A2", "HON107", "S_016", 45666.44991898148),
    @("D_NEHAMER001", "C_PRIYATIR001", "7943FA0E-BA47-4710-880A-606294BFC73D", "MER105", "S_017", 45674.58295138889),
    @("D_BALAJITAT001", "C_LAKSHMIERO001", "658DF915-E21B-4057-AB64-7A948ADEE559", "TAT106", "S_018", 45760.42690972222),
    @("D_DEEPAKMAR001", "C_PREETIVEL001", "F82E7F95-C5DD-470C-BD78-616874E1D889", "MAR109", "S_019", 45700.56508101852),
    @("D_SNEHAMAH001", "C_DEEPATIR001", "02CC8F43-A3DA-48B8-AE4C-DE9AAFA8350F", "MAH104", "S_020", 45778.46548611111),
    @("D_BALAJITAT001", "C_SANDEEPTIR001", "7878650D-C981-46A1-B5AD-8BD971B5273B", "TAT106", "S_021", 45696.4675462963),
    @("D_DEEPAKMAR001", "C_RAHULSAL001", "385E7458-D624-4B08-ADB7-159826CD8BF5", "MAR109", "S_022", 45693.56962962963),
    @("D_POOJAMAR001", "C_VIVEKTHA001", "FD52A531-5E2B-4F6A-BD6A-920F49BA8F34", "MAR109", "S_023", 45682.70148148148),
    @("D_SANJAYHYU001", "C_MEENAKAN001", "CA1121F0-6D7E-40F9-9362-6D5D726DE29B", "HYU101", "S_024", 45726.60721064815),
    @("D_RAJESHHYU001", "C_AMITABHERO001", "EA27201E-64A7-4625-987E-CDCA9CC1914B", "HYU101", "S_025", 45768.67391203704),
    @("D_SANJAYHYU001", "C_PREETIVEL001", "052A0B28-647F-4F28-BB69-1B19F435D8DC", "HYU101", "S_026", 45793.73270833334),
    @("D_VIDYAMAR001", "C_SANDEEPTIR001", "08E3F3BA-8377-4E49-9F71-2B50BCDE07F7", "MAR109", "S_027", 45693.63803240741),
    @("D_SIVAKUMARMAH001", "C_VIVEKTHA001", "22574029-F597-4F35-B4E8-B16D7C23A86A", "MAH104", "S_028", 45710.41773148148),
    @("D_BALAJITAT001", "C_SWETHACOI001", "282C509F-537C-4CFD-8DC5-E8EDC5D76FCA", "TAT106", "S_029", 45703.60927083333),
    @("D_KARTHIKFOR001", "C_DEEPATIR001", "1F2615B2-C7BB-4A58-A884-B19D53F6A597", "FOR108", "S_030", 45743.59858796297)
)

$dateFormat = $ws.Cells.Item(2, 6).NumberFormat()

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $dateCell = $ws.Cells.Item($r, 6)
    $dateCell.Value = $row[5]
    $dateCell.NumberFormat = $dateFormat
    $r = $r + 1
}
